$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp in A1 (COVID data refresh)
$ws.Range("A1").Value = "Datos actualizados a 18 de Octubre de 2020 a las 12:44"

# Refresh the per-country statistics (Casos totales, Nuevos casos,
# Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
# for every row whose source numbers changed in this data refresh.
# Row/country pairings are unchanged - only the B:H figures move.

# Row 20: Banglades
$ws.Range("B20").Value = 388569
$ws.Range("C20").Value = 1274
$ws.Range("D20").Value = 303972
$ws.Range("E20").Value = 78937
$ws.Range("G20").Value = 14
$ws.Range("H20").Value = 5660

# Row 27: Israel
$ws.Range("B27").Value = 302911
$ws.Range("C27").Value = 141
$ws.Range("D27").Value = 267221
$ws.Range("E27").Value = 33488
$ws.Range("G27").Value = 12
$ws.Range("H27").Value = 2202

# Row 32: Rumania
$ws.Range("B32").Value = 180388
$ws.Range("C32").Value = 3920
$ws.Range("D32").Value = 130894
$ws.Range("E32").Value = 43622
$ws.Range("G32").Value = 60
$ws.Range("H32").Value = 5872

# Row 38: Nepal
$ws.Range("B38").Value = 129431
$ws.Range("C38").Value = 204
$ws.Range("D38").Value = 126406
$ws.Range("E38").Value = 2801
$ws.Range("G38").Value = 1
$ws.Range("H38").Value = 224

# Row 39: Catar
$ws.Range("B39").Value = 129304
$ws.Range("D39").Value = 89840
$ws.Range("E39").Value = 38737
$ws.Range("H39").Value = 727

# Row 42: Kuwait
$ws.Range("B42").Value = 115602
$ws.Range("C42").Value = 1215
$ws.Range("D42").Value = 107516
$ws.Range("E42").Value = 7623
$ws.Range("G42").Value = 4
$ws.Range("H42").Value = 463

# Row 43: Emiratos Arabes Unidos
$ws.Range("B43").Value = 115483
$ws.Range("D43").Value = 107108
$ws.Range("E43").Value = 7681
$ws.Range("H43").Value = 694

# Row 44: Kazajistan
$ws.Range("B44").Value = 109953
$ws.Range("C44").Value = 1657
$ws.Range("D44").Value = 95624
$ws.Range("E44").Value = 13228
$ws.Range("G44").Value = 30
$ws.Range("H44").Value = 1101

# Row 45: Oman
$ws.Range("B45").Value = 109406
$ws.Range("C45").Value = 104
$ws.Range("D45").Value = 105001
$ws.Range("E45").Value = 2637
$ws.Range("H45").Value = 1768

# Row 69: Irlanda
$ws.Range("B69").Value = 48790
$ws.Range("C69").Value = 945
$ws.Range("D69").Value = 26889
$ws.Range("E69").Value = 21176
$ws.Range("G69").Value = 26
$ws.Range("H69").Value = 725

# Row 70: Libia
$ws.Range("B70").Value = 48678
$ws.Range("D70").Value = 23364
$ws.Range("E70").Value = 23465
$ws.Range("H70").Value = 1849

# Row 92: Costa de Marfil
$ws.Range("B92").Value = 20498
$ws.Range("C92").Value = 871
$ws.Range("D92").Value = 13262
$ws.Range("E92").Value = 7049
$ws.Range("G92").Value = 7
$ws.Range("H92").Value = 187

# Row 93: Malasia
$ws.Range("B93").Value = 20301
$ws.Range("D93").Value = 19983
$ws.Range("E93").Value = 197
$ws.Range("H93").Value = 121

# Row 100: Senegal
$ws.Range("B100").Value = 15418
$ws.Range("C100").Value = 26
$ws.Range("D100").Value = 13814
$ws.Range("E100").Value = 1287

# Row 103: Eslovenia
$ws.Range("B103").Value = 13142
$ws.Range("C103").Value = 726
$ws.Range("D103").Value = 6313
$ws.Range("E103").Value = 6641
$ws.Range("G103").Value = 4
$ws.Range("H103").Value = 188

# Row 110: Uganda
$ws.Range("B110").Value = 10590
$ws.Range("C110").Value = 135
$ws.Range("D110").Value = 6992
$ws.Range("E110").Value = 3501
$ws.Range("G110").Value = 1
$ws.Range("H110").Value = 97

# Row 126: Sri Lanka
$ws.Range("D126").Value = 3403
$ws.Range("E126").Value = 2059

# Row 176: Gibraltar
$ws.Range("B176").Value = 571
$ws.Range("C176").Value = 13
$ws.Range("D176").Value = 455
$ws.Range("E176").Value = 116

Write-Host "Update complete"
